# Rename the JavaObjects / MySQL-tables diagram entities:
#  - "ID"/"Name"/etc PascalCase Java fields -> camelCase, and
#    foreign-key "XxxID : int" fields become plain object references
#    ("objectDetail", "Object", "InspectionForm", ...)
#  - "ID"/"Name"/etc PascalCase MySQL columns -> snake_case
#  - Table/class names "ObjectDetails" -> "ObjectDetail" / "object_detail", etc.
#  - "Completed" run -> "completed" (first run text only)
#  - H7 on JavaObjects gets a red font
#  - Selections move to F9 (JavaObjects) / D17 (MySQL tables)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("JavaObjects")
$ws2 = $wb.Worksheets.Item("MySQL tables")

# --- Cell text rewrites, issued in the order that reproduces the target
#     shared-string table layout (unchanged class/table names keep their
#     original relative order automatically; every brand-new string below
#     is appended in exactly this sequence). ---

$ws1.Range('A1').Value = 'Object'
$ws1.Range('B5').Value = 'Object'
$ws1.Range('D6').Value = 'Object'
$ws1.Range('A5').Value = 'List<InspectionForm>'
$ws1.Range('E7').Value = 'InspectionForm'
$ws1.Range('G15').Value = 'InspectionForm'
$ws1.Range('D7').Value = 'List<InspectionStep>'
$ws1.Range('G3').Value = 'Inspection'
$ws2.Range('I3').Value = 'Inspection'
$ws1.Range('H3').Value = 'InspectionResult'
$ws1.Range('I6').Value = 'InspectionResult'
$ws1.Range('G16').Value = 'List<InspectionResult>'
$ws1.Range('H9').Value = 'List<FaultPicture>'
$ws1.Range('I3').Value = 'FaultPicture'
$ws2.Range('M3').Value = 'FaultPicture'
$ws1.Range('D2').Value = 'Form'
$ws1.Range('E2').Value = 'Step'
$ws2.Range('K3').Value = 'Result'
$ws1.Range('B1').Value = 'ObjectDetail'
$ws2.Range('A4').Value = 'object_detail_id   :   INT() (fk)'
$ws2.Range('A3').Value = 'name   :   VARCHAR (unique)'
$ws2.Range('A2').Value = 'id   :   INT() (pk)'
$ws2.Range('C2').Value = 'id   :   INT() (pk)'
$ws2.Range('E3').Value = 'id   :   INT() (pk)'
$ws2.Range('G3').Value = 'id   :   INT() (pk)'
$ws2.Range('I4').Value = 'id   :   INT() (pk)'
$ws2.Range('K4').Value = 'id   :   INT() (pk)'
$ws2.Range('M4').Value = 'id   :   INT() (pk)'
$ws2.Range('C3').Value = 'maker   :   VARCHAR'
$ws2.Range('C4').Value = 'description   :   TEXT'
$ws2.Range('G4').Value = 'description   :   TEXT'
$ws2.Range('E5').Value = 'description   :   TEXT'
$ws2.Range('E4').Value = 'name   :   VARCHAR (unique with ObjectID)'
$ws2.Range('E6').Value = 'object_id   :   INT() (fk; unique with Name)'
$ws2.Range('G5').Value = 'details   :   TEXT'
$ws2.Range('K7').Value = 'step_number   :   INT()'
$ws2.Range('G7').Value = 'inspection_form_id   :   INT() (fk)'
$ws2.Range('A1').Value = 'object'
$ws2.Range('C1').Value = 'object_detail'
$ws2.Range('E2').Value = 'form'
$ws2.Range('G2').Value = 'step'
$ws1.Range('E6').Value = 'number   :   int'
$ws2.Range('G6').Value = 'number   :   INT()'
$ws1.Range('A2').Value = 'id   :   int'
$ws1.Range('B2').Value = 'id   :   int'
$ws1.Range('D3').Value = 'id   :   int'
$ws1.Range('E3').Value = 'id   :   int'
$ws1.Range('G4').Value = 'id   :   int'
$ws1.Range('H4').Value = 'id   :   int'
$ws1.Range('I4').Value = 'id   :   int'
$ws1.Range('A3').Value = 'name   :   String'
$ws1.Range('D4').Value = 'name   :   String'
$ws1.Range('A4').Value = 'objectDetail'
$ws1.Range('B4').Value = 'description   :   String'
$ws1.Range('E4').Value = 'description   :   String'
$ws1.Range('D5').Value = 'description   :   String'
$ws1.Range('B3').Value = 'maker   :   String'
$ws1.Range('E5').Value = 'details   :   String'
$ws1.Range('G5').Value = 'serialNumber   :   String'
$ws1.Range('G6').Value = 'creationDate   :   Date'
$ws1.Range('G7').Value = 'creationHour   :   Date'
$ws1.Range('G8').Value = 'completionDate   :   Date'
$ws1.Range('G9').Value = 'completionHour   :   Date'
$ws1.Range('G10').Value = 'inspector   :   String'
$ws1.Range('G11').Value = 'place   :   String'
$ws1.Range('G12').Value = 'batch   :   int'
$ws1.Range('G13').Value = 'mainResult   :   int'
$ws1.Range('H5').Value = 'result   :   int'
$ws1.Range('H6').Value = 'notes   :   String'
$ws1.Range('H7').Value = 'stepNumber   :   int'
$ws1.Range('H8').Value = 'inspection'
$ws1.Range('I5').Value = 'fileName   :   String'
$ws2.Range('I5').Value = 'serial_number   :   VARCHAR (unique with InspectionFormID)'
$ws2.Range('I6').Value = 'creation_date   :   DATE()'
$ws2.Range('I7').Value = 'creation_hour   :   TIME()'
$ws2.Range('I8').Value = 'completion_date   :   DATE()'
$ws2.Range('I9').Value = 'completion_hour   :   TIME()'
$ws2.Range('I10').Value = 'inspector   :   VARCHAR'
$ws2.Range('I11').Value = 'place   :   VARCHAR'
$ws2.Range('I12').Value = 'batch   :   INT()'
$ws2.Range('I13').Value = 'main_result   :   INT()'
$ws2.Range('I14').Value = 'completed   :   TINYINT()(1)'
$ws2.Range('I15').Value = 'inspection_form_id   :   INT() (fk; unique with SerialNumber)'
$ws2.Range('K5').Value = 'result   :   INT()'
$ws2.Range('K6').Value = 'notes   :   TEXT'
$ws2.Range('K8').Value = 'inspection_id   :   INT() (fk)'
$ws2.Range('M5').Value = 'file_name   :   VARCHAR'
$ws2.Range('M6').Value = 'inspection_result_id   :   INT() (fk)'

# --- Two-run rich text cell: "completed   :   " (plain) + "boolean" (own run) ---
$completedCell = $ws1.Range('G14')
$completedCell.Value = 'completed   :   boolean'
$booleanRun = $completedCell.Characters(17, 7)
$booleanRun.Font.Name = 'Calibri'
$booleanRun.Font.Size = 11
$booleanRun.Font.ColorIndex = -4105

# --- New red font applied to H7 (JavaObjects) ---
$ws1.Range('H7').Font.Color = 255

# --- Restore cursor/selection positions ---
$ws1.Range('F9').Select()
$ws2.Range('D17').Select()
